$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 60500
$ws.Range("J3").Value = 60500
$ws.Range("L3").Value = 60500
$ws.Range("N3").Value = -60728
$ws.Range("H18").Value = 331.25
$ws.Range("I18").Value = 331.25
$ws.Range("K18").Value = 331.25
$ws.Range("M18").Value = -47.25
$ws.Range("H64").Value = 1800
$ws.Range("J64").Value = 1800
$ws.Range("L64").Value = 1800
$ws.Range("N64").Value = -2296
$ws.Range("H67").Value = 1800
$ws.Range("J67").Value = 1800
$ws.Range("L67").Value = 1800
$ws.Range("N67").Value = -3516
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H102").Value = 60500
$ws.Range("J102").Value = 60500
$ws.Range("L102").Value = 60500
$ws.Range("N102").Value = -66990
$ws.Range("H116").Value = 3987.5
$ws.Range("I116").Value = 3987.5
$ws.Range("K116").Value = 3987.5
$ws.Range("M116").Value = -545.5
$ws.Range("H132").Value = 1078.2941
$ws.Range("I132").Value = 1095.3125
$ws.Range("K132").Value = 3285.9375
$ws.Range("M132").Value = -755.9375
$ws.Range("H138").Value = 1469
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("H141").Value = 190
$ws.Range("I141").Value = 190
$ws.Range("K141").Value = 570
$ws.Range("M141").Value = 4610

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 918.2
$ws.Range("I32").Value = 918.2
$ws.Range("K32").Value = 918.2
$ws.Range("M32").Value = -631.2
$ws.Range("H74").Value = 251895.88
$ws.Range("I74").Value = 251895.88
$ws.Range("K74").Value = 251895.88
$ws.Range("M74").Value = -251021.88
$ws.Range("H77").Value = 251895.88
$ws.Range("I77").Value = 251895.88
$ws.Range("K77").Value = 1259479.4
$ws.Range("M77").Value = -1255111.4
$ws.Range("H80").Value = 79501
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H82").Value = 34499.5
$ws.Range("J82").Value = 48999
$ws.Range("L82").Value = 48999
$ws.Range("N82").Value = -49721
$ws.Range("H83").Value = 79501
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H85").Value = 34499.5
$ws.Range("J85").Value = 48999
$ws.Range("L85").Value = 48999
$ws.Range("N85").Value = -51495
$ws.Range("H96").Value = 37092
$ws.Range("J96").Value = 37092
$ws.Range("L96").Value = 37092
$ws.Range("N96").Value = -42584
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H122").Value = 8917
$ws.Range("J122").Value = 7432.1665
$ws.Range("L122").Value = 22296.4995
$ws.Range("N122").Value = -27196.4995

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("H94").Value = 2454.4546
$ws.Range("I94").Value = 1999.8334
$ws.Range("J94").Value = 3000
$ws.Range("K94").Value = 1999.8334
$ws.Range("L94").Value = 3000
$ws.Range("M94").Value = -1548.8334
$ws.Range("N94").Value = -3902

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 15845.167
$ws.Range("I38").Value = 2769
$ws.Range("K38").Value = 2769
$ws.Range("M38").Value = -2392
$ws.Range("H39").Value = 27878.857
$ws.Range("I39").Value = 11050.667
$ws.Range("K39").Value = 11050.667
$ws.Range("M39").Value = -10659.667
$ws.Range("H41").Value = 35000
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 35000
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 35000
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -35856
$ws.Range("H46").Value = 15845.167
$ws.Range("I46").Value = 2769
$ws.Range("K46").Value = 2769
$ws.Range("M46").Value = -2558
$ws.Range("H49").Value = 27878.857
$ws.Range("I49").Value = 11050.667
$ws.Range("K49").Value = 11050.667
$ws.Range("M49").Value = -10868.667
$ws.Range("H68").Value = 78509
$ws.Range("J68").Value = 75636.25
$ws.Range("L68").Value = 75636.25
$ws.Range("N68").Value = -77134.25
$ws.Range("H71").Value = 78509
$ws.Range("J71").Value = 75636.25
$ws.Range("L71").Value = 226908.75
$ws.Range("N71").Value = -234396.75
$ws.Range("H95").Value = 36712.43
$ws.Range("J95").Value = 36712.43
$ws.Range("L95").Value = 36712.43
$ws.Range("N95").Value = -42204.43
$ws.Range("H107").Value = 915.6667
$ws.Range("I107").Value = 748.75
$ws.Range("K107").Value = 748.75
$ws.Range("M107").Value = 1171.25
$ws.Range("H132").Value = 1550.3334
$ws.Range("I132").Value = 1465.0588
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 4395.1764
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -1865.1764
$ws.Range("N132").Value = -14060
$ws.Range("H134").Value = 977.6667
$ws.Range("I134").Value = 977.6667
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 2933.0001
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -398.0001000000002
$ws.Range("N134").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 750
$ws.Range("J107").Value = 820
$ws.Range("L107").Value = 2460
$ws.Range("N107").Value = -6300

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
$ws.Range("H59").Value = 43999.2
$ws.Range("J59").Value = 43999.2
$ws.Range("L59").Value = 43999.2
$ws.Range("N59").Value = -45165.2
$ws.Range("H62").Value = 90000
$ws.Range("I62").Value = 90000
$ws.Range("K62").Value = 90000
$ws.Range("M62").Value = -89314
$ws.Range("H65").Value = 90000
$ws.Range("I65").Value = 90000
$ws.Range("K65").Value = 270000
$ws.Range("M65").Value = -266568
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10000
$ws.Range("I7").Value = 10400
$ws.Range("K7").Value = 10400
$ws.Range("M7").Value = -10288
$ws.Range("H20").Value = 9500
$ws.Range("J20").Value = 9500
$ws.Range("L20").Value = 9500
$ws.Range("N20").Value = -9952
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H126").Value = 10000
$ws.Range("I126").Value = 10400
$ws.Range("K126").Value = 31200
$ws.Range("M126").Value = -28730
$ws.Range("H132").Value = 2506.4666
$ws.Range("I132").Value = 2874.625
$ws.Range("J132").Value = 2085.7144
$ws.Range("K132").Value = 8623.875
$ws.Range("L132").Value = 6257.1432
$ws.Range("M132").Value = -6093.875
$ws.Range("N132").Value = -11317.1432

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 40000
$ws.Range("J110").Value = 40000
$ws.Range("L110").Value = 40000
$ws.Range("N110").Value = -48180
$ws.Range("H126").Value = 2857.3333
$ws.Range("I126").Value = 2857.3333
$ws.Range("K126").Value = 8571.999899999999
$ws.Range("M126").Value = -6101.999899999999
$ws.Range("H132").Value = 2229.9688
$ws.Range("I132").Value = 1702.1072
$ws.Range("K132").Value = 5106.321599999999
$ws.Range("M132").Value = -2576.321599999999
